# Weekly update for "Hortaliza, Vega Monumental Concepción - Cilantro":
# a new week's pair of observations (Primera/Segunda) is recorded at the
# top of the data block (rows 62-63), every existing observation shifts
# down by one week (two rows), and the oldest pair that falls off the
# bottom (old rows 156-157) is appended as the two new last rows
# (158-159).
#
# Only columns D (Fecha) and O (Origen) actually change value per row;
# every other column is constant across all observations in this sheet,
# so re-stamping D/O after duplicating the last pair reproduces the
# target state exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 62
$lastDataRow = 157
$newLastRow = 159

# 1) Grow the table: duplicate the current last two rows (156-157) into
#    the two brand-new rows (158-159) so all formatting/styles/values
#    for the other columns come along for free.
$ws.Range("A156:R156").Copy($ws.Range("A158:R158"))
$ws.Range("A157:R157").Copy($ws.Range("A159:R159"))

# 2) Snapshot the original Fecha (D) / Origen (O) columns for the data
#    block before mutating anything.
$dVals = @{}
$oVals = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $dVals[$r] = $ws.Cells.Item($r, 4).Value2
    $oVals[$r] = $ws.Cells.Item($r, 15).Value2
}

# 3) Shift every existing pair down by two rows: new row r gets the
#    Fecha/Origen that used to live at row r-2. Walk from the bottom up
#    so each write's source has already been captured in step 2.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 2); $r--) {
    $ws.Cells.Item($r, 4).Value = $dVals[$r - 2]
    $ws.Cells.Item($r, 15).Value = $oVals[$r - 2]
}

# 4) Stamp the brand-new top pair (rows 62-63) with this week's values.
$ws.Cells.Item($firstDataRow, 4).Value = 44579
$ws.Cells.Item($firstDataRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($firstDataRow + 1, 4).Value = 44579
$ws.Cells.Item($firstDataRow + 1, 15).Value = "Región Metropolitana"
